$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2L4LT4ST_37T")

# Fill in the torque/index data for columns C:F, rows 2-17 (ordered row-by-row
# to mirror the commit's sheetData order).
$data = @(
    @(2,  0.09, 0,    0.17, 0),
    @(3,  0.19, 0.36, 0.17, 0.335),
    @(4,  0.19, 0.17, 0.21, 0.04),
    @(5,  0.09, 0.29, 0.14, 0.41),
    @(6,  0.15, 0.07, 0,    0.16),
    @(7,  0.11, 0.17, 0.12, 0.24),
    @(8,  0.18, 0.16, 0.12, 0.13),
    @(9,  0.17, 0.18, 0.21, 0),
    @(10, 0.12, 0.02, 0.19, 0.16),
    @(11, 0.1,  0.17, 0.11, 0.15),
    @(12, 0.27, 0.15, 0.1,  0.11),
    @(13, 0.18, 0.14, 0.21, 0.27),
    @(14, 0,    0.18, 0.23, 0),
    @(15, 0.13, 0.12, 0.13, 0.27),
    @(16, 0.23, 0.11, 0,    0),
    @(17, 0.23, 0.15, 0.16, 0.75)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[1]
    $ws.Cells.Item($row, 4).Value = $entry[2]
    $ws.Cells.Item($row, 5).Value = $entry[3]
    $ws.Cells.Item($row, 6).Value = $entry[4]
}

# Update the selection on the sheet
$ws.Range("M11").Select()

$wb.Save()
